# Developer Guide diagram fix: correct the border-cut-off issue by nudging
# the diagram group + one connector, and split the editBond() signature
# text into annotated runs (one run per identifier / literal segment).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Top level diagram group ("Group 107") ---------------------------------
# a:off x="98646" y="22069"  ->  x="53165" y="21266"
$g = $s.Shapes.Item(1)
$g.Left = 4.186225
$g.Top  = 1.6745

# --- "Curved Connector 31" inside the group ---------------------------------
# a:off x="5596798" y="2305559"  ->  x="5596798" y="2312539" (x unchanged)
$cc = $g.GroupItems.Item(28)
$cc.Top = 182.08968503937007

# --- "TextBox 40" (the investmentEditBond(...) signature) ------------------
$tb = $g.GroupItems.Item(37)
$tr = $tb.TextFrame.TextRange

# Paragraph 2 currently reads "(bondName, year, rate, ui)" as a single run.
# Re-split it into: "(" / "bondName" / ", year, rate, " / "ui" / ")"
$para2 = $tr.Paragraphs(2)
$para2.Characters(1, 1).Text  = "("
$para2.Characters(2, 8).Text  = "bondName"
$para2.Characters(10, 14).Text = ", year, rate, "
$para2.Characters(24, 2).Text = "ui"
$para2.Characters(26, 1).Text = ")"

# Re-splitting the run nudges the shape's auto-fit height slightly; restore
# the original box height so only the text runs (not the geometry) changed.
$tb.Height = 34.203307086614174
